# Update the crypto price/volume table (cryptos list refresh, per
# GitHub Actions commit "Updated cryptos list ... with GitHub Actions").
#
# Note: several "Price" column values (D) are plain numeric-looking
# strings (e.g. "580.67", "1.00", "0.0751"). The source workbook stores
# these as text, so we prefix them with a leading apostrophe to force
# Excel to keep them as text instead of silently converting them to
# numbers (exactly what typing `'580.67` into a cell does in the UI).
# Values that already contain two dots (e.g. "64.567.44") or spaces
# (the Volume(1h) percentages) are never auto-converted by Excel, so
# they are assigned as plain strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "64.567.44"
$ws.Range("E2").Value = "  -2.07%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "3.417.14"
$ws.Range("E3").Value = "  -3.15%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.07%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "'580.67"
$ws.Range("E5").Value = "  -3.81%  "

# --- Row 6: Solana ---
$ws.Range("D6").Value = "'134.20"
$ws.Range("E6").Value = "  -6.40%  "

# --- Row 7: USDC ---
$ws.Range("D7").Value = "'1.00"

# --- Row 8: LidoStakedEther ---
$ws.Range("D8").Value = "3.416.44"
$ws.Range("E8").Value = "  -3.08%  "

# --- Row 9: XRP ---
$ws.Range("E9").Value = "  -5.46%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  -8.66%  "

# --- Row 11: Toncoin ---
$ws.Range("D11").Value = "'6.97"
$ws.Range("E11").Value = "  -10.91%  "

# --- Row 12: Cardano ---
$ws.Range("E12").Value = "  -8.06%  "

# --- Row 13: WrappedliquidstakedEther2.0 ---
$ws.Range("D13").Value = "3.999.12"
$ws.Range("E13").Value = "  -3.31%  "

# --- Row 14: ShibaInu ---
$ws.Range("E14").Value = "  -8.46%  "

# --- Row 15: WrappedEther ---
$ws.Range("D15").Value = "3.431.30"
$ws.Range("E15").Value = "  -2.34%  "

# --- Row 17: Avalanche ---
$ws.Range("D17").Value = "'25.88"
$ws.Range("E17").Value = "  -8.45%  "

# --- Row 18: WrappedBTC ---
$ws.Range("D18").Value = "64.545.08"
$ws.Range("E18").Value = "  -1.92%  "

# --- Row 19: Uniswap ---
$ws.Range("D19").Value = "'9.38"
$ws.Range("E19").Value = "  -14.07%  "

# --- Row 20: Polkadot ---
$ws.Range("D20").Value = "'5.71"
$ws.Range("E20").Value = "  -7.61%  "

# --- Row 21: Chainlink ---
$ws.Range("D21").Value = "'13.39"
$ws.Range("E21").Value = "  -7.80%  "

# --- Row 22: BitcoinCash ---
$ws.Range("D22").Value = "'377.41"
$ws.Range("E22").Value = "  -10.09%  "

# --- Row 24: Polygon ---
$ws.Range("E24").Value = "  -9.36%  "

# --- Row 25: Litecoin ---
$ws.Range("D25").Value = "'71.33"
$ws.Range("E25").Value = "  -6.91%  "

# --- Row 26: WrappedeETH ---
$ws.Range("D26").Value = "3.556.35"
$ws.Range("E26").Value = "  -2.91%  "

# --- Row 27: PEPE ---
$ws.Range("E27").Value = "  -8.15%  "

# --- Row 29: RenderToken ---
$ws.Range("E29").Value = "  -8.81%  "

# --- Row 30/31: InternetComputer(DFINITY) and PancakeSwap swap places ---
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.17"
$ws.Range("E30").Value = "  -11.84%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'7.93"
$ws.Range("E31").Value = "  -10.65%  "

# --- Row 32: RenzoRestakedETH ---
$ws.Range("D32").Value = "3.435.50"
$ws.Range("E32").Value = "  -2.76%  "

# --- Row 33: USDe ---
$ws.Range("E33").Value = "  +0.02%  "

# --- Row 34: EthereumClassic ---
$ws.Range("E34").Value = "  -5.26%  "

# --- Row 35: Kaspa ---
$ws.Range("E35").Value = "  -9.38%  "

# --- Row 36: Monero ---
$ws.Range("D36").Value = "'168.29"
$ws.Range("E36").Value = "  -4.13%  "

# --- Row 37: Fetch.AI ---
$ws.Range("E37").Value = "  -13.52%  "

# --- Row 38: Aptos ---
$ws.Range("D38").Value = "'6.63"
$ws.Range("E38").Value = "  -11.61%  "

# --- Row 39: ImmutableX ---
$ws.Range("E39").Value = "  -10.60%  "

# --- Row 40: NEARProtocol ---
$ws.Range("D40").Value = "'4.54"
$ws.Range("E40").Value = "  -13.05%  "

# --- Row 41: Hedera ---
$ws.Range("D41").Value = "'0.0751"
$ws.Range("E41").Value = "  -7.62%  "

# --- Row 42: Mantle ---
$ws.Range("D42").Value = "'0.804"
$ws.Range("E42").Value = "  -6.07%  "

# --- Row 43: FirstDigitalUSD ---
$ws.Range("E43").Value = "  +0.09%  "

# --- Row 44: OKB ---
$ws.Range("D44").Value = "'41.55"
$ws.Range("E44").Value = "  -8.49%  "

# --- Row 45: Filecoin ---
$ws.Range("E45").Value = "  -14.30%  "

# --- Row 46: Stacks ---
$ws.Range("E46").Value = "  -9.80%  "

# --- Row 47: ONDO ---
$ws.Range("D47").Value = "'1.10"
$ws.Range("E47").Value = "  -1.06%  "

# --- Row 48: EnergySwap ---
$ws.Range("D48").Value = "'22.38"
$ws.Range("E48").Value = "  -3.23%  "

# --- Row 49: Cosmos ---
$ws.Range("D49").Value = "'6.42"
$ws.Range("E49").Value = "  -8.46%  "

# --- Row 50: Maker ---
$ws.Range("D50").Value = "2.165.42"
$ws.Range("E50").Value = "  -6.58%  "

# --- Row 51: dogwifhat ---
$ws.Range("E51").Value = "  -16.74%  "
